# ---------------------------------------------------------------------------
# Menu() is not working properly
#
# "items" sheet: unit/quantity corrections + drop the two trailing rows that
# are no longer stocked (pen drive / Ups).
#
# "soldProduct" sheet: three new sales were recorded (Pen, Computer, Pen
# drive) -> append as rows 22-24.
#
# "userAccount" sheet: a new user (Usrah saba) registered -> append as row 9.
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

function Set-TextValue {
    param($range, [string]$value)
    # Force a numeric-looking string to be stored as text (matches the
    # source file's inlineStr cells, e.g. phone numbers / passwords), then
    # drop back to the Normal cell style so no stray formatting lingers.
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

# ---------------------------------------------------------------------------
# items
# ---------------------------------------------------------------------------
$wsItems = $wb.Worksheets.Item("items")

$wsItems.Range("D2").Value = 198
$wsItems.Range("D3").Value = 9
$wsItems.Range("D5").Value = 10
$wsItems.Range("C8").Value = 110
$wsItems.Range("D8").Value = 10
$wsItems.Range("D9").Value = 10
$wsItems.Range("D12").Value = 10
$wsItems.Range("C13").Value = 10
$wsItems.Range("D13").Value = 10000

# Drop rows 14 (pen drive) and 15 (Ups) - delete bottom-up so indices stay valid
$wsItems.Rows.Item(15).Delete()
$wsItems.Rows.Item(14).Delete()

# ---------------------------------------------------------------------------
# soldProduct
# ---------------------------------------------------------------------------
$wsSold = $wb.Worksheets.Item("soldProduct")

$wsSold.Range("A22").Value = 1000
$wsSold.Range("B22").Value = "Pen"
$wsSold.Range("C22").Value = 2
$wsSold.Range("D22").Value = 10
$wsSold.Range("E22").Value = 20
$wsSold.Range("F22").Value = "Sadia afroz"
$wsSold.Range("G22").Value = "fsdfsd"
Set-TextValue $wsSold.Range("H22") "5435"
$wsSold.Range("I22").Value = "nogod"

$wsSold.Range("A23").Value = 1001
$wsSold.Range("B23").Value = "Computer"
$wsSold.Range("C23").Value = 1
$wsSold.Range("D23").Value = 10000
$wsSold.Range("E23").Value = 10000
$wsSold.Range("F23").Value = "Sadia"
$wsSold.Range("G23").Value = "uttara"
Set-TextValue $wsSold.Range("H23") "5435325"
$wsSold.Range("I23").Value = "COD"

$wsSold.Range("A24").Value = 1019
$wsSold.Range("B24").Value = "Pen drive"
$wsSold.Range("C24").Value = 2
$wsSold.Range("D24").Value = 1000
$wsSold.Range("E24").Value = 2000
$wsSold.Range("F24").Value = "Usrah saba"
$wsSold.Range("G24").Value = "dflsaj"
Set-TextValue $wsSold.Range("H24") "543534"
$wsSold.Range("I24").Value = "bkash"

# ---------------------------------------------------------------------------
# userAccount
# ---------------------------------------------------------------------------
$wsUsers = $wb.Worksheets.Item("userAccount")

$wsUsers.Range("A9").Value = "Usrah saba"
$wsUsers.Range("B9").Value = "usrah"
Set-TextValue $wsUsers.Range("C9") "123"
$wsUsers.Range("D9").Value = "gfdt"
Set-TextValue $wsUsers.Range("E9") "4534"
